# Generate Report for Handback
#
# The "1b18accd-ef67-4ed0-b431-a21bf8f620ba.md" localization job has been
# handed back (its translations are now in sync with en-US). Reflect that
# on the Overview sheet and on each per-locale sheet (zh-cn, de-de):
#   - Status flips from "Ready for handoff" to "Handed back: in sync with en-US"
#   - The "Latest Target File" / "Latest Handback File" columns get populated
#     (as hyperlinks, matching the look of the existing "Latest Handoff File"
#     column) and "Latest Handback DateTime" gets stamped.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: row 2 is the 1b18accd-...md file.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus

# ---------------------------------------------------------------------
# Helper: update one locale sheet (zh-cn / de-de) + its handback timestamp.
# ---------------------------------------------------------------------
function Update-LocaleSheet([string]$sheetName, [string]$handbackStamp) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Status column (B) for the 1b18accd-...md row.
    $ws.Range("B2").Value = $newStatus

    # Grab the existing source-file / handoff-file hyperlink targets so the
    # new "target" / "handback" hyperlinks point at the same place, the way
    # the other rows in this sheet do.
    $srcLink = $null
    $xlfLink = $null
    foreach ($h in $ws.Hyperlinks) {
        $addr = $h.Range.Address()
        if ($addr -eq "$A$2") { $srcLink = $h }
        if ($addr -eq "$C$2") { $xlfLink = $h }
    }

    $srcDisplay = $ws.Range("A2").Value2
    $xlfDisplay = $ws.Range("C2").Value2
    $srcAddress = $srcLink.Address
    $xlfAddress = $xlfLink.Address

    # Latest Target File (E2) - mirrors the Source File Name hyperlink.
    $ws.Hyperlinks.Add($ws.Range("E2"), $srcAddress, [Type]::Missing, [Type]::Missing, $srcDisplay) | Out-Null

    # Latest Handback File (F2) - mirrors the Latest Handoff File hyperlink.
    $ws.Hyperlinks.Add($ws.Range("F2"), $xlfAddress, [Type]::Missing, [Type]::Missing, $xlfDisplay) | Out-Null

    # Latest Handback DateTime (G2).
    $ws.Range("G2").Value = $handbackStamp
}

Update-LocaleSheet "zh-cn" "2016-03-09 00:48:23"
Update-LocaleSheet "de-de" "2016-03-09 00:48:54"
